$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 19.80936873547465
$ws.Range("D2").Value = 3.568290086743094
$ws.Range("E2").Value = 29.69072338448435
$ws.Range("F2").Value = 19.30527164261118
$ws.Range("G2").Value = 3.588097124384748
$ws.Range("L2").Value = 12.42501740481856
$ws.Range("O2").Value = 17.05282932392036

# Row 3
$ws.Range("B3").Value = 19.26206573308494
$ws.Range("D3").Value = 3.533650201644278
$ws.Range("E3").Value = 28.59691468878408
$ws.Range("F3").Value = 19.3196031058285
$ws.Range("G3").Value = 3.590844547852015
$ws.Range("L3").Value = 11.92829353618152
$ws.Range("O3").Value = 17.15160235732488

# Row 4
$ws.Range("B4").Value = 18.91857585666316
$ws.Range("D4").Value = 3.512500154624606
$ws.Range("E4").Value = 27.9079351199038
$ws.Range("F4").Value = 19.34006569217391
$ws.Range("G4").Value = 3.592617534756709
$ws.Range("L4").Value = 11.61153301778787
$ws.Range("O4").Value = 17.22137557870469

# Row 5
$ws.Range("B5").Value = 18.77690396961934
$ws.Range("D5").Value = 3.503918660245035
$ws.Range("E5").Value = 27.62325643173564
$ws.Range("F5").Value = 19.35131250189764
$ws.Range("G5").Value = 3.593361757011936
$ws.Range("L5").Value = 11.47963179487107
$ws.Range("O5").Value = 17.2520766542692

# Row 6
$ws.Range("B6").Value = 18.75328230804909
$ws.Range("D6").Value = 3.502496179422128
$ws.Range("E6").Value = 27.57576359453945
$ws.Range("F6").Value = 19.35335487235422
$ws.Range("G6").Value = 3.593486648348285
$ws.Range("L6").Value = 11.45756393325662
$ws.Range("O6").Value = 17.25731072047542

# Row 7
$ws.Range("B7").Value = 18.91667185148832
$ws.Range("D7").Value = 3.512384260751328
$ws.Range("E7").Value = 27.90411105447494
$ws.Range("F7").Value = 19.34020563126931
$ws.Range("G7").Value = 3.592627483573557
$ws.Range("L7").Value = 11.60976536500907
$ws.Range("O7").Value = 17.22178047764169

# Row 8
$ws.Range("B8").Value = 19.62230472810172
$ws.Range("D8").Value = 3.556326006481325
$ws.Range("E8").Value = 29.31744357267058
$ws.Range("F8").Value = 19.30777892870134
$ws.Range("G8").Value = 3.589026620468004
$ws.Range("L8").Value = 12.25626885970665
$ws.Range("O8").Value = 17.08497615652931

# Row 9
$ws.Range("B9").Value = 20.93992693783237
$ws.Range("D9").Value = 3.643114656589415
$ws.Range("E9").Value = 31.93343193812035
$ws.Range("F9").Value = 19.33761981838838
$ws.Range("G9").Value = 3.582644693788773
$ws.Range("L9").Value = 13.42536940863237
$ws.Range("O9").Value = 16.89027952557721

# Row 10
$ws.Range("B10").Value = 21.8591008359723
$ws.Range("D10").Value = 3.706822701170229
$ws.Range("E10").Value = 33.74000895413678
$ws.Range("F10").Value = 19.41739612810009
$ws.Range("G10").Value = 3.578365155256472
$ws.Range("L10").Value = 14.21816878375835
$ws.Range("O10").Value = 16.79361472127692

# Row 11
$ws.Range("B11").Value = 22.26507764915804
$ws.Range("D11").Value = 3.735702458092447
$ws.Range("E11").Value = 34.53329131447788
$ws.Range("F11").Value = 19.46634204532838
$ws.Range("G11").Value = 3.576506102808075
$ws.Range("L11").Value = 14.56349263014026
$ws.Range("O11").Value = 16.76001333753799

# Row 12
$ws.Range("B12").Value = 22.41694748478939
$ws.Range("D12").Value = 3.746616434545603
$ws.Range("E12").Value = 34.82933905372424
$ws.Range("F12").Value = 19.4866964474061
$ws.Range("G12").Value = 3.575814663661171
$ws.Range("L12").Value = 14.69198823055836
$ws.Range("O12").Value = 16.74880381653908

# Row 13
$ws.Range("B13").Value = 22.3843242624571
$ws.Range("D13").Value = 3.744267019156004
$ws.Range("E13").Value = 34.76577697810927
$ws.Range("F13").Value = 19.48223188879772
$ws.Range("G13").Value = 3.575963020712857
$ws.Range("L13").Value = 14.66441637204375
$ws.Range("O13").Value = 16.7511502732725

# Row 14
$ws.Range("B14").Value = 22.27761006936207
$ws.Range("D14").Value = 3.736600844535775
$ws.Range("E14").Value = 34.55773584097846
$ws.Range("F14").Value = 19.46798016438839
$ws.Range("G14").Value = 3.576448966711191
$ws.Range("L14").Value = 14.57410990464652
$ws.Range("O14").Value = 16.75906064166819

# Row 15
$ws.Range("B15").Value = 22.21199844487976
$ws.Range("D15").Value = 3.731901982742746
$ws.Range("E15").Value = 34.42973122699254
$ws.Range("F15").Value = 19.45948745374835
$ws.Range("G15").Value = 3.576748254090726
$ws.Range("L15").Value = 14.51849700715821
$ws.Range("O15").Value = 16.76410387631887

# Row 16
$ws.Range("B16").Value = 21.83231441030393
$ws.Range("D16").Value = 3.704932649698792
$ws.Range("E16").Value = 33.68756991403149
$ws.Range("F16").Value = 19.41445215257453
$ws.Range("G16").Value = 3.578488407923831
$ws.Range("L16").Value = 14.19528668550824
$ws.Range("O16").Value = 16.79602137067597

# Row 17
$ws.Range("B17").Value = 21.59618972779076
$ws.Range("D17").Value = 3.688356433784028
$ws.Range("E17").Value = 33.22478188067024
$ws.Range("F17").Value = 19.39006722729384
$ws.Range("G17").Value = 3.579578354697939
$ws.Range("L17").Value = 13.99303126967554
$ws.Range("O17").Value = 16.81827513310767

# Row 18
$ws.Range("B18").Value = 21.45923868969885
$ws.Range("D18").Value = 3.678813138201886
$ws.Range("E18").Value = 32.95592478897822
$ws.Range("F18").Value = 19.37723316336388
$ws.Range("G18").Value = 3.580213525209847
$ws.Range("L18").Value = 13.87526019002454
$ws.Range("O18").Value = 16.83204986865226

# Row 19
$ws.Range("B19").Value = 21.41267763612082
$ws.Range("D19").Value = 3.675580613889919
$ws.Range("E19").Value = 32.86444328125152
$ws.Range("F19").Value = 19.37309233490993
$ws.Range("G19").Value = 3.58043000422531
$ws.Range("L19").Value = 13.83513998537443
$ws.Range("O19").Value = 16.83688047903524

# Row 20
$ws.Range("B20").Value = 21.62144429827523
$ws.Range("D20").Value = 3.690121997002212
$ws.Range("E20").Value = 33.27432506714696
$ws.Range("F20").Value = 19.39253971421764
$ws.Range("G20").Value = 3.579461473447535
$ws.Range("L20").Value = 14.01471116836439
$ws.Range("O20").Value = 16.81580510551154

# Row 21
$ws.Range("B21").Value = 22.30900609294507
$ws.Range("D21").Value = 3.738853246969117
$ws.Range("E21").Value = 34.6189624595447
$ws.Range("F21").Value = 19.47211688252623
$ws.Range("G21").Value = 3.576305892732996
$ws.Range("L21").Value = 14.60069722414514
$ws.Range("O21").Value = 16.75669589222289

# Row 22
$ws.Range("B22").Value = 22.74745225209211
$ws.Range("D22").Value = 3.770569363079272
$ws.Range("E22").Value = 35.47230731074177
$ws.Range("F22").Value = 19.53472788553469
$ws.Range("G22").Value = 3.574316619192155
$ws.Range("L22").Value = 14.97041143435527
$ws.Range("O22").Value = 16.72690245507497

# Row 23
$ws.Range("B23").Value = 22.51447946639619
$ws.Range("D23").Value = 3.753656470092796
$ws.Range("E23").Value = 35.01926338209692
$ws.Range("F23").Value = 19.50034239728674
$ws.Range("G23").Value = 3.575371668729784
$ws.Range("L23").Value = 14.7743210616531
$ws.Range("O23").Value = 16.7419878639483

# Row 24
$ws.Range("B24").Value = 21.61003043511211
$ws.Range("D24").Value = 3.689323827037402
$ws.Range("E24").Value = 33.25193527641314
$ws.Range("F24").Value = 19.39141821130281
$ws.Range("G24").Value = 3.579514288836242
$ws.Range("L24").Value = 14.00491432745182
$ws.Range("O24").Value = 16.81691875180574

# Row 25
$ws.Range("B25").Value = 20.59148768527281
$ws.Range("D25").Value = 3.61961466059787
$ws.Range("E25").Value = 31.24467807219973
$ws.Range("F25").Value = 19.31943411114393
$ws.Range("G25").Value = 3.584298950800937
$ws.Range("L25").Value = 13.12033619587173
$ws.Range("O25").Value = 16.93491181605881
